# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.244.92"
$ws.Range("E2").Value = "  +3.22%  "

$ws.Range("D3").Value = "3.134.05"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'220.09"
$ws.Range("E5").Value = "  +4.56%  "

$ws.Range("D6").Value = "'622.63"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "'0.378"
$ws.Range("E7").Value = "  +1.20%  "

$ws.Range("D8").Value = "'0.910"
$ws.Range("E8").Value = "  +10.16%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "3.130.37"
$ws.Range("E10").Value = "  +1.54%  "

$ws.Range("D11").Value = "'0.746"
$ws.Range("E11").Value = "  +25.55%  "

$ws.Range("D12").Value = "'0.190"
$ws.Range("E12").Value = "  +6.48%  "

$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "  +5.95%  "

$ws.Range("D14").Value = "'34.10"
$ws.Range("E14").Value = "  +7.56%  "

$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").Value = "91.085.20"
$ws.Range("E16").Value = "  +3.46%  "

$ws.Range("D17").Value = "3.701.77"
$ws.Range("E17").Value = "  +1.47%  "

$ws.Range("D18").Value = "3.142.20"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("D19").Value = "'3.83"
$ws.Range("E19").Value = "  +19.13%  "

$ws.Range("D20").Value = "'0.0000226"
$ws.Range("E20").Value = "  +6.01%  "

$ws.Range("D21").Value = "'14.15"
$ws.Range("E21").Value = "  +7.77%  "

$ws.Range("D22").Value = "'432.66"
$ws.Range("E22").Value = "  +2.53%  "

$ws.Range("E23").Value = "  +7.66%  "

$ws.Range("D24").Value = "'5.17"
$ws.Range("E24").Value = "  +6.54%  "

$ws.Range("E25").Value = "  +11.84%  "

$ws.Range("D26").Value = "'12.40"
$ws.Range("E26").Value = "  +8.44%  "

$ws.Range("D27").Value = "'83.73"
$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("D28").Value = "3.281.19"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").Value = "'0.167"
$ws.Range("E30").Value = "  +7.73%  "

$ws.Range("D31").Value = "'9.06"
$ws.Range("E31").Value = "  +12.57%  "

$ws.Range("E32").Value = "  -16.15%  "

$ws.Range("D33").Value = "'528.82"
$ws.Range("E33").Value = "  +4.71%  "

$ws.Range("D34").Value = "'3.85"
$ws.Range("E34").Value = "  +8.17%  "

$ws.Range("D35").Value = "'7.38"
$ws.Range("E35").Value = "  +11.44%  "

$ws.Range("E36").Value = "  +8.95%  "

$ws.Range("E37").Value = "  +6.52%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'23.42"
$ws.Range("E38").Value = "  +5.19%  "

$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "'1.86"
$ws.Range("E39").Value = "  +2.60%  "

$ws.Range("D40").Value = "'22.30"
$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0773"
$ws.Range("E42").Value = "  +17.54%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'0.144"
$ws.Range("E44").Value = "  +6.05%  "

$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").Value = "'0.378"
$ws.Range("E45").Value = "  +5.55%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.92"
$ws.Range("E46").Value = "  +5.55%  "

$ws.Range("D47").Value = "'143.79"
$ws.Range("E47").Value = "  -3.41%  "

$ws.Range("D48").Value = "'44.19"
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("E49").Value = "  +10.76%  "

$ws.Range("D50").Value = "'0.000264"
$ws.Range("E50").Value = "  +22.96%  "

$ws.Range("D51").Value = "'168.19"
$ws.Range("E51").Value = "  +7.63%  "
